$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two duplicate rows (original rows 13 and 10).
# Delete the lower one first so the row number of the other deletion
# target doesn't shift.
$ws.Rows("13").Delete()
$ws.Rows("10").Delete()

# Edit remaining fields in row 9 (Luis Jonathan Diaz Mattus)
$ws.Range("B9").Value = "Luis con apellido"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "10000000000"

$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "2.023"

$ws.Range("N9").Value = 6

$ws.Range("O9").Value = "Enseñanza aa"
